$d = $word.ActiveDocument

# 1) Remove the old _GoBack bookmark (after " De Silva" on the title page).
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 2) Fix the "NPL" -> "NLP" typo and split the run, inserting a fresh
#    _GoBack bookmark right at the edit point (mirrors Word's own
#    behaviour of dropping a _GoBack bookmark at the last edited spot).
#    The bookmark is inserted FIRST (splitting the run surgically) so the
#    later, narrower Find/Replace only rewrites the small isolated run
#    instead of the whole paragraph.
$r = $d.Content
$r.Find.Execute("Natural Language Processing (NPL", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$splitPoint = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $splitPoint) | Out-Null

$fix = $d.Content
$fix.Find.Execute("NPL", $true, $false, $false, $false, $false, $true, 1, $false, "NLP", 2) | Out-Null
